$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leetcode")

# Append a new row (31) with a new "Binary Search" problem entry.
# (New shared strings are appended in the order the values are first
# written, so "Fun one." must be set before "Koko Eating Bananas" to match
# the target shared-string table ordering.)
$ws.Range("A31").Value = "Leetcode"
$ws.Range("B31").Value = 875
$ws.Range("D31").Value = "Binary Search"
$ws.Range("E31").Value = "Medium"
$ws.Range("F31").Value = "Neetcode 150"
$ws.Range("G31").Value = "SOLVED"
$ws.Range("H31").Value = "14/06/2025"
$ws.Range("I31").Value = "Fun one."
$ws.Range("C31").Value = "Koko Eating Bananas"

# Match the styling used by the rest of the table (left align identifier,
# wrap text on notes column).
$ws.Range("B31").HorizontalAlignment = -4131
$ws.Range("I31").WrapText = $true

# Move the active selection to the newly added row / reset horizontal scroll.
$ws.Range("D31").Select()
